$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 188
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()
$ws.Range("H36").Value = 188
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()
$ws.Range("H40").Value = 2800.3333
$ws.Range("J40").Value = 3200
$ws.Range("L40").Value = 3200
$ws.Range("N40").Value = -3550
$ws.Range("H46").Value = 3999
$ws.Range("I46").Value = 3999
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 11997
$ws.Range("L46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -11878
$ws.Range("H48").Value = 2345
$ws.Range("I48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("M48").ClearContents()
$ws.Range("H50").Value = 500
$ws.Range("I50").Value = 500
$ws.Range("K50").Value = 1500
$ws.Range("M50").Value = -1025
$ws.Range("H56").Value = 2345
$ws.Range("I56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("M56").ClearContents()
$ws.Range("H60").Value = 3999
$ws.Range("I60").Value = 3999
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 11997
$ws.Range("L60").Value = 0
$ws.Range("M60").ClearContents()
$ws.Range("N60").Value = -11513
$ws.Range("H116").Value = 2150
$ws.Range("I116").Value = 2150
$ws.Range("K116").Value = 2150
$ws.Range("M116").Value = 1292
$ws.Range("H138").Value = 4187.45
$ws.Range("J138").Value = 4085.2354
$ws.Range("L138").Value = 12255.7062
$ws.Range("N138").Value = -22535.7062

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").ClearContents()
$ws.Range("N49").Value = 0
$ws.Range("H50").Value = 43800
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 43800
$ws.Range("K50").Value = 0
$ws.Range("L50").ClearContents()
$ws.Range("M50").Value = 43800
$ws.Range("N50").Value = -45228
$ws.Range("H63").Value = 1374.1666
$ws.Range("I63").Value = 712.5
$ws.Range("J63").Value = 2697.5
$ws.Range("K63").Value = 712.5
$ws.Range("L63").Value = 2697.5
$ws.Range("M63").Value = -26.5
$ws.Range("N63").Value = -4069.5
$ws.Range("H66").Value = 1374.1666
$ws.Range("I66").Value = 712.5
$ws.Range("J66").Value = 2697.5
$ws.Range("K66").Value = 3562.5
$ws.Range("L66").Value = 13487.5
$ws.Range("M66").Value = -130.5
$ws.Range("N66").Value = -20351.5
$ws.Range("H122").Value = 3493.5
$ws.Range("I122").Value = 2991.3333
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 8973.999899999999
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -6523.999899999999
$ws.Range("N122").Value = -19900
$ws.Range("H132").Value = 1968
$ws.Range("I132").Value = 1986
$ws.Range("J132").Value = 1938
$ws.Range("K132").Value = 5958
$ws.Range("L132").Value = 5814
$ws.Range("M132").Value = -3428
$ws.Range("N132").Value = -10874

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1899.091
$ws.Range("I99").Value = 1812.8572
$ws.Range("K99").Value = 1812.8572
$ws.Range("M99").Value = -314.8571999999999
$ws.Range("H105").Value = 3165.8147
$ws.Range("I105").Value = 2936.5417
$ws.Range("K105").Value = 2936.5417
$ws.Range("M105").Value = -1189.5417

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 483.83334
$ws.Range("I25").Value = 465
$ws.Range("J25").Value = 502.66666
$ws.Range("K25").Value = 465
$ws.Range("L25").Value = 502.66666
$ws.Range("M25").Value = -291
$ws.Range("N25").Value = -850.66666
$ws.Range("H35").Value = 1175
$ws.Range("I35").Value = 1175
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 1175
$ws.Range("L35").Value = 0
$ws.Range("M35").ClearContents()
$ws.Range("N35").Value = -881
$ws.Range("H36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("L36").ClearContents()
$ws.Range("N36").Value = 0
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").ClearContents()
$ws.Range("N40").Value = 0
$ws.Range("H55").Value = 7333.3335
$ws.Range("J55").Value = 10000
$ws.Range("L55").Value = 10000
$ws.Range("N55").Value = -10630
$ws.Range("H58").Value = 1231
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()
$ws.Range("H134").Value = 2101.3333
$ws.Range("I134").Value = 1902.25
$ws.Range("K134").Value = 5706.75
$ws.Range("M134").Value = -3171.75
$ws.Range("H136").Value = 1231
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 1001.5
$ws.Range("J12").Value = 1673.375
$ws.Range("L12").Value = 5020.125
$ws.Range("N12").Value = -5366.125
$ws.Range("H33").Value = 300.5
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()
$ws.Range("H50").Value = 408.5
$ws.Range("I50").Value = 189
$ws.Range("J50").Value = 1506
$ws.Range("K50").Value = 567
$ws.Range("L50").Value = 4518
$ws.Range("M50").Value = -86
$ws.Range("N50").Value = -5480
$ws.Range("H53").Value = 408.5
$ws.Range("I53").Value = 189
$ws.Range("J53").Value = 1506
$ws.Range("K53").Value = 567
$ws.Range("L53").Value = 4518
$ws.Range("M53").Value = -86
$ws.Range("N53").Value = -5480
$ws.Range("H117").Value = 532.3333
$ws.Range("I117").Value = 498.5
$ws.Range("J117").Value = 600
$ws.Range("K117").Value = 1495.5
$ws.Range("L117").Value = 1800
$ws.Range("M117").Value = 1946.5
$ws.Range("N117").Value = -8684

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H54").Value = 15000
$ws.Range("J54").Value = 15000
$ws.Range("L54").Value = 15000
$ws.Range("N54").Value = -15780

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 317.16666
$ws.Range("I16").Value = 317.16666
$ws.Range("K16").Value = 317.16666
$ws.Range("M16").Value = -147.16666
$ws.Range("H29").Value = 22500
$ws.Range("I29").Value = 22500
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 22500
$ws.Range("L29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("N29").Value = -22205
$ws.Range("H31").Value = 3575.7144
$ws.Range("I31").Value = 3575.7144
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 3575.7144
$ws.Range("L31").Value = 0
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -3327.7144
$ws.Range("H46").Value = 2121.923
$ws.Range("I46").Value = 1561.0416
$ws.Range("J46").Value = 3019.3333
$ws.Range("K46").Value = 1561.0416
$ws.Range("L46").Value = 3019.3333
$ws.Range("M46").Value = -1373.0416
$ws.Range("N46").Value = -3395.3333
$ws.Range("H47").Value = 65
$ws.Range("J47").Value = 65
$ws.Range("L47").Value = 65
$ws.Range("N47").Value = -1045
$ws.Range("H52").Value = 65
$ws.Range("J52").Value = 65
$ws.Range("L52").Value = 65
$ws.Range("N52").Value = -531
$ws.Range("H55").Value = 688
$ws.Range("J55").Value = 759.06665
$ws.Range("L55").Value = 759.06665
$ws.Range("N55").Value = -1105.06665
$ws.Range("H58").Value = 9634
$ws.Range("I58").Value = 7899.5
$ws.Range("K58").Value = 7899.5
$ws.Range("M58").Value = -7639.5
$ws.Range("H136").Value = 3994
$ws.Range("I136").Value = 3356.6667
$ws.Range("J136").Value = 4950
$ws.Range("K136").Value = 10070.0001
$ws.Range("L136").Value = 14850
$ws.Range("M136").Value = -7520.000100000001
$ws.Range("N136").Value = -19950

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 714.2857
$ws.Range("I81").Value = 333.66666
$ws.Range("J81").Value = 999.75
$ws.Range("K81").Value = 667.33332
$ws.Range("L81").Value = 1999.5
$ws.Range("M81").Value = 393.66668
$ws.Range("N81").Value = -4121.5
$ws.Range("H84").Value = 714.2857
$ws.Range("I84").Value = 333.66666
$ws.Range("J84").Value = 999.75
$ws.Range("K84").Value = 3336.6666
$ws.Range("L84").Value = 9997.5
$ws.Range("M84").Value = 1967.3334
$ws.Range("N84").Value = -20605.5
$ws.Range("H136").Value = 1076.238
$ws.Range("I136").Value = 1094.8422
$ws.Range("J136").Value = 899.5
$ws.Range("K136").Value = 3284.5266
$ws.Range("L136").Value = 2698.5
$ws.Range("M136").Value = -734.5266000000001
$ws.Range("N136").Value = -7798.5
